$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B20").Value = 6164
$ws.Range("D20").Value = 5568491
$ws.Range("E20").Value = 903.3891953277093
$ws.Range("F20").Value = 6.477802729314219
$ws.Range("H20").Value = 25.96477957089338
